$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actividades durante la pasantia")

# New activities appended to the calendar, plus the fix to row 76's date
# (previously a literal text "23(07/2022" stored in column A; now a real
# date value like every other row).
$rows = @(
    @{ Row = 76; Serial = 44765; Activity = "Reunion con vero organizando el codigo con datos vacios"; Hours = 1.5 },
    @{ Row = 77; Serial = 44768; Activity = "Apoyo pruebas suiza, curso manejo de datos, reunion con vero del articulo"; Hours = 6 },
    @{ Row = 78; Serial = 44775; Activity = "Reunion con vero y luisa, curso manejo de datos"; Hours = 3 },
    @{ Row = 79; Serial = 44777; Activity = "curso manejo de datos"; Hours = 2 },
    @{ Row = 80; Serial = 44780; Activity = "Revision articulos reactividad"; Hours = 2 },
    @{ Row = 81; Serial = 44782; Activity = "Reunion con vero y luisa reactividad, curso manejo de datos"; Hours = 4 },
    @{ Row = 82; Serial = 44783; Activity = "Parque explora apropiacion social del conocimiento"; Hours = 3 },
    @{ Row = 83; Serial = 44784; Activity = "Finalizacion del curso de manejo de datos con herramientas de googles"; Hours = 2 },
    @{ Row = 84; Serial = 44789; Activity = "Graficos ICC solo por bandas"; Hours = 2 }
)

# Column A in row 75 already carries the date number format/style we want
# (m/d/yyyy) - copy it down onto each new date cell so they reuse the same
# style id instead of Excel minting a new number-format/style combo.
foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A75").Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false

    $ws.Cells.Item($row, 1).Value = $r.Serial
    $ws.Cells.Item($row, 2).Value = $r.Activity
    $ws.Cells.Item($row, 3).Value = $r.Hours
}

$ws.Range("A73").Select()
$ws.Range("C84").Select()
